# Append the new resale-numbers row (2024-01-23 21:29:05) as row 88.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 88

# Columns A:D hold text-like values ("2024-01-23", "21:29:05", "Tuesday",
# "03") that Excel would otherwise auto-coerce into a date serial / time
# serial / number (losing the leading zero on "03"). Force them to be
# stored as literal text, matching every other row in the sheet, then
# drop the temporary format override so no stray cell style is left
# behind (existing data rows carry no explicit style).
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-23"
$ws.Cells.Item($row, 2).Value = "21:29:05"
$ws.Cells.Item($row, 3).Value = "Tuesday"
$ws.Cells.Item($row, 4).Value = "03"

$textRange.ClearFormats()

# Columns E:T are the per-city numeric resale counts.
$ws.Cells.Item($row, 5).Value = 138484
$ws.Cells.Item($row, 6).Value = 141295
$ws.Cells.Item($row, 7).Value = 171147
$ws.Cells.Item($row, 8).Value = 148885
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 123371
$ws.Cells.Item($row, 11).Value = 223736
$ws.Cells.Item($row, 12).Value = 256245
$ws.Cells.Item($row, 13).Value = 185038
$ws.Cells.Item($row, 14).Value = 110211
$ws.Cells.Item($row, 15).Value = 41335
$ws.Cells.Item($row, 16).Value = 30896
$ws.Cells.Item($row, 17).Value = 73466
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42647
$ws.Cells.Item($row, 20).Value = -1
